# ISW_2021_1C_LinksClasesGrabadas.xlsx
# "Se actualiza link de clases grabadas"
#
# Updates the "Segundo Parcial" schedule block (rows 28-31):
#   - Row 28 (was "Comparacion de enfoques tradicional, lean y agile")
#     becomes the "Retrospectiva" class, with an extra instruction
#     sentence and a new "Clase Grabada Retrospectiva" recorded-class
#     link in column F. The row grows taller to fit the wrapped text.
#   - Row 29 ("Publicidad en Instagram...") gains a new recorded-class
#     link "Clase Grabada Practico 13" in column G.
#   - Row 31 ("Revisiones Tecnicas") gains a new recorded-class link
#     "Clase Grabada Revisiones Tecnicas" in column G.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Hyperlink colour used throughout the sheet for "clase grabada" links
$hyperlinkColor = 13391121   # RGB(0x11,0x55,0xCC) -> FF1155CC

# ---------------------------------------------------------------
# Row 28 - Retrospectiva
# ---------------------------------------------------------------
$ws.Range("D28").Value2 = "Retrospectiva"

$ws.Range("E28").Value2 = "Deben traer para trabajar en clase un cuadro comparativo de los 3 enfoques (ver TP 14 en la guía de prácticos). Es evaluable. Se usará como recuperatorio de los TP conceptuales, para los grupos que lo necesiten"
$ws.Range("E28").WrapText = $true

$ws.Range("F28").Value2 = ""
$linkF28 = $ws.Hyperlinks.Add($ws.Range("F28"), "https://youtu.be/5Kq1DxN0q1o", "", "", "Clase Grabada Retrospectiva")
$ws.Range("F28").Font.Underline = 2
$ws.Range("F28").Font.Color = $hyperlinkColor

# Taller row to fit the wrapped instructions text
$ws.Rows.Item(28).RowHeight = 54

# ---------------------------------------------------------------
# Row 29 - new recorded-class link (Practico 13)
# ---------------------------------------------------------------
$ws.Range("G29").Value2 = ""
$linkG29 = $ws.Hyperlinks.Add($ws.Range("G29"), "https://youtu.be/3k4sGz1b2qY", "", "", "Clase Grabada Practico 13")
$ws.Range("G29").Font.Underline = 2
$ws.Range("G29").Font.Color = $hyperlinkColor
$ws.Range("G29").HorizontalAlignment = -4108
$ws.Range("G29").VerticalAlignment = -4108

# ---------------------------------------------------------------
# Row 31 - new recorded-class link (Revisiones Tecnicas)
# ---------------------------------------------------------------
$ws.Range("G31").Value2 = ""
$linkG31 = $ws.Hyperlinks.Add($ws.Range("G31"), "https://youtu.be/7hE9vW2mYxA", "", "", "Clase Grabada Revisiones Tecnicas")
$ws.Range("G31").Font.Underline = 2
$ws.Range("G31").Font.Color = $hyperlinkColor
$ws.Range("G31").HorizontalAlignment = -4108
$ws.Range("G31").VerticalAlignment = -4108
